$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.05878958822046783, 0.6848687919290307, 1.944388015115169, 1.394413143625364, 1.397936283545429, 147)
    3 = @(0.06771567608964886, 0.6998484907537829, 1.910406263097688, 1.382174469123811, 1.385266920986894, 146)
    4 = @(0.07543715453382718, 0.6915666599019875, 1.885412560133963, 1.373103259093781, 1.375781748473243, 145)
    5 = @(0.07888146508037401, 0.7033649776424507, 1.96991753687704,  1.403537508183176, 1.406210292228227, 144)
    6 = @(0.08557990013155226, 0.712304576464059,  1.96830456032932,  1.402962779381306, 1.405272345080091, 143)
    7 = @(0.0938298966656947,  0.6980419742238877, 1.929693431924385, 1.389134058298329, 1.390867611014577, 142)
    8 = @(0.1082541482841633,  0.6906349894607012, 1.935468717112932, 1.391211241010125, 1.391937795370696, 141)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
